$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C ("Förändrad") from 45182 to 45184 for all existing data rows (2..410)
for ($r = 2; $r -le 410; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value = 45184
    }
}

# Ensure row 410's explicit row height is persisted (matches original sheet's row metadata)
$ws.Rows.Item(410).RowHeight = 15

# 2. Append the new row 411 with its data
$ws.Cells.Item(411, 1).Value = "A 43037-2023"
$ws.Cells.Item(411, 2).Value = 45182
$ws.Cells.Item(411, 3).Value = 45184
$ws.Cells.Item(411, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(411, 5).Value = "VARBERG"
$ws.Cells.Item(411, 7).Value = 3.2
$ws.Cells.Item(411, 8).Value = 0
$ws.Cells.Item(411, 9).Value = 0
$ws.Cells.Item(411, 10).Value = 0
$ws.Cells.Item(411, 11).Value = 0
$ws.Cells.Item(411, 12).Value = 0
$ws.Cells.Item(411, 13).Value = 0
$ws.Cells.Item(411, 14).Value = 0
$ws.Cells.Item(411, 15).Value = 0
$ws.Cells.Item(411, 16).Value = 0
$ws.Cells.Item(411, 17).Value = 0

# Match date formatting / styling used by the Datum (B) and Förändrad (C) columns in prior rows
$ws.Cells.Item(411, 2).NumberFormat = $ws.Cells.Item(410, 2).NumberFormat
$ws.Cells.Item(411, 3).NumberFormat = $ws.Cells.Item(410, 3).NumberFormat

# Match the wrap-text style used by column R (Artnamn) in prior rows, leaving it blank
$ws.Cells.Item(411, 18).WrapText = $true
